# Update the dSF column (F) values per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -4
    4  = -2
    5  = 2
    7  = -1
    9  = -3
    10 = -1
    11 = 5
    12 = -2
    13 = -1
    14 = -3
    16 = -4
    17 = -3
    18 = -2
    19 = 1
    20 = 3
    21 = 1
    22 = 5
    23 = -1
    24 = 4
    25 = 6
    26 = 13
    27 = 4
    28 = 6
    29 = 1
    30 = -1
    31 = -1
    32 = -4
    33 = -4
    34 = -2
    35 = 4
    36 = -3
    37 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
